$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New day's competitor data (2025-07-04, serial 45842) appended below existing rows.
$newRowsData = @(
  @(45842, "МилЛимон",    92,     5,     9),
  @(45842, "Фоксфорд",    268719, 31044, 92910),
  @(45842, "Умскул",      492489, 68726, 529767),
  @(45842, "Skysmart",    70072,  18337, 370211),
  @(45842, "ЕГЭLAND",     41253,  14076, 3092),
  @(45842, "Алгоритмика", 62898,  35449, 45812),
  @(45842, "easycode",    27781,  17930, 1793),
  @(45842, "Котокод",     4055,   254,   16)
)

$startRow = 50
for ($i = 0; $i -lt $newRowsData.Count; $i++) {
  $r = $startRow + $i
  $row = $newRowsData[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
}

# Copy the formatting from the previous last data row (49) down onto the
# freshly added rows so they pick up the existing date / text number formats
# instead of minting brand new style entries.
$ws.Range("A49:E49").Copy()
$ws.Range("A50:E57").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Bold the header row.
$ws.Range("A1:H1").Font.Bold = $true

# Restore scroll position / selection that Excel persists on save.
$ws.Range("B3").Select() | Out-Null
